$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 9).Value = 'aa'
$ws.Cells.Item(2, 10).Value = 'Agree/Accept'
$ws.Cells.Item(6, 9).Value = 'sd'
$ws.Cells.Item(6, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(8, 9).Value = 'sv'
$ws.Cells.Item(8, 10).Value = 'Statement-opinion'
$ws.Cells.Item(45, 9).Value = 'sd'
$ws.Cells.Item(45, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(54, 9).Value = 'b'
$ws.Cells.Item(54, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(81, 9).Value = 'sd'
$ws.Cells.Item(81, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(102, 9).Value = 'b'
$ws.Cells.Item(102, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(104, 9).Value = 'b'
$ws.Cells.Item(104, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(109, 9).Value = 'b'
$ws.Cells.Item(109, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(110, 9).Value = 'sd'
$ws.Cells.Item(110, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(113, 9).Value = 'b'
$ws.Cells.Item(113, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(132, 9).Value = 'sv'
$ws.Cells.Item(132, 10).Value = 'Statement-opinion'
$ws.Cells.Item(137, 9).Value = 'sd'
$ws.Cells.Item(137, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(147, 9).Value = 'sd'
$ws.Cells.Item(147, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(149, 9).Value = 'b'
$ws.Cells.Item(149, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(156, 9).Value = 'sd'
$ws.Cells.Item(156, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(167, 9).Value = 'sv'
$ws.Cells.Item(167, 10).Value = 'Statement-opinion'
$ws.Cells.Item(168, 9).Value = '%'
$ws.Cells.Item(168, 10).Value = 'Uninterpretable'
$ws.Cells.Item(169, 9).Value = 'b'
$ws.Cells.Item(169, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(172, 9).Value = '%'
$ws.Cells.Item(172, 10).Value = 'Uninterpretable'
$ws.Cells.Item(173, 9).Value = 'sd'
$ws.Cells.Item(173, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(189, 9).Value = 'b'
$ws.Cells.Item(189, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(191, 9).Value = 'sv'
$ws.Cells.Item(191, 10).Value = 'Statement-opinion'
$ws.Cells.Item(194, 9).Value = 'sd'
$ws.Cells.Item(194, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(205, 9).Value = 'sv'
$ws.Cells.Item(205, 10).Value = 'Statement-opinion'
$ws.Cells.Item(206, 9).Value = 'ba'
$ws.Cells.Item(206, 10).Value = 'Appreciation'
$ws.Cells.Item(215, 9).Value = 'sd'
$ws.Cells.Item(215, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(222, 9).Value = 'sd'
$ws.Cells.Item(222, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(224, 9).Value = 'sd'
$ws.Cells.Item(224, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(228, 9).Value = 'aa'
$ws.Cells.Item(228, 10).Value = 'Agree/Accept'
$ws.Cells.Item(229, 9).Value = 'sd'
$ws.Cells.Item(229, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(239, 9).Value = 'ba'
$ws.Cells.Item(239, 10).Value = 'Appreciation'
$ws.Cells.Item(266, 9).Value = 'sd'
$ws.Cells.Item(266, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(269, 9).Value = 'b'
$ws.Cells.Item(269, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(279, 9).Value = 'ba'
$ws.Cells.Item(279, 10).Value = 'Appreciation'
$ws.Cells.Item(288, 9).Value = 'sd'
$ws.Cells.Item(288, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(303, 9).Value = 'sd'
$ws.Cells.Item(303, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(311, 9).Value = 'sd'
$ws.Cells.Item(311, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(315, 9).Value = 'sd'
$ws.Cells.Item(315, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(331, 9).Value = 'sd'
$ws.Cells.Item(331, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(336, 9).Value = 'sd'
$ws.Cells.Item(336, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(340, 9).Value = 'sd'
$ws.Cells.Item(340, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(353, 9).Value = 'sv'
$ws.Cells.Item(353, 10).Value = 'Statement-opinion'
$ws.Cells.Item(361, 9).Value = 'b'
$ws.Cells.Item(361, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(376, 9).Value = 'ba'
$ws.Cells.Item(376, 10).Value = 'Appreciation'
$ws.Cells.Item(379, 9).Value = 'sv'
$ws.Cells.Item(379, 10).Value = 'Statement-opinion'
$ws.Cells.Item(380, 9).Value = 'sd'
$ws.Cells.Item(380, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(381, 9).Value = 'sv'
$ws.Cells.Item(381, 10).Value = 'Statement-opinion'
$ws.Cells.Item(393, 9).Value = 'sd'
$ws.Cells.Item(393, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(395, 9).Value = 'sv'
$ws.Cells.Item(395, 10).Value = 'Statement-opinion'
$ws.Cells.Item(399, 9).Value = 'ba'
$ws.Cells.Item(399, 10).Value = 'Appreciation'
$ws.Cells.Item(416, 9).Value = 'sv'
$ws.Cells.Item(416, 10).Value = 'Statement-opinion'
$ws.Cells.Item(421, 9).Value = 'ba'
$ws.Cells.Item(421, 10).Value = 'Appreciation'
$ws.Cells.Item(435, 9).Value = 'sd'
$ws.Cells.Item(435, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(444, 9).Value = 'b'
$ws.Cells.Item(444, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(454, 9).Value = 'sd'
$ws.Cells.Item(454, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(458, 9).Value = 'b'
$ws.Cells.Item(458, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(463, 9).Value = 'b'
$ws.Cells.Item(463, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(473, 9).Value = 'sd'
$ws.Cells.Item(473, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(474, 9).Value = 'sd'
$ws.Cells.Item(474, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(477, 9).Value = '%'
$ws.Cells.Item(477, 10).Value = 'Uninterpretable'
$ws.Cells.Item(480, 9).Value = 'b'
$ws.Cells.Item(480, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(481, 9).Value = 'b'
$ws.Cells.Item(481, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(482, 9).Value = 'b'
$ws.Cells.Item(482, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(483, 9).Value = 'qy'
$ws.Cells.Item(483, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(487, 9).Value = 'sd'
$ws.Cells.Item(487, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(488, 9).Value = 'qy'
$ws.Cells.Item(488, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(490, 9).Value = 'sd'
$ws.Cells.Item(490, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(493, 9).Value = 'sd'
$ws.Cells.Item(493, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(500, 9).Value = 'b'
$ws.Cells.Item(500, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(515, 9).Value = 'ba'
$ws.Cells.Item(515, 10).Value = 'Appreciation'
$ws.Cells.Item(547, 9).Value = 'sv'
$ws.Cells.Item(547, 10).Value = 'Statement-opinion'
$ws.Cells.Item(548, 9).Value = 'sd'
$ws.Cells.Item(548, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(569, 9).Value = 'sd'
$ws.Cells.Item(569, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(571, 9).Value = 'sd'
$ws.Cells.Item(571, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(573, 9).Value = 'aa'
$ws.Cells.Item(573, 10).Value = 'Agree/Accept'
$ws.Cells.Item(574, 9).Value = 'aa'
$ws.Cells.Item(574, 10).Value = 'Agree/Accept'
$ws.Cells.Item(578, 9).Value = 'b'
$ws.Cells.Item(578, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(580, 9).Value = 'b'
$ws.Cells.Item(580, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(596, 9).Value = 'aa'
$ws.Cells.Item(596, 10).Value = 'Agree/Accept'
$ws.Cells.Item(611, 9).Value = 'sd'
$ws.Cells.Item(611, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(617, 9).Value = 'sv'
$ws.Cells.Item(617, 10).Value = 'Statement-opinion'
$ws.Cells.Item(624, 9).Value = 'b'
$ws.Cells.Item(624, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(632, 9).Value = 'sv'
$ws.Cells.Item(632, 10).Value = 'Statement-opinion'
$ws.Cells.Item(642, 9).Value = 'aa'
$ws.Cells.Item(642, 10).Value = 'Agree/Accept'
$ws.Cells.Item(651, 9).Value = 'sv'
$ws.Cells.Item(651, 10).Value = 'Statement-opinion'
$ws.Cells.Item(653, 9).Value = 'b'
$ws.Cells.Item(653, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(672, 9).Value = 'sd'
$ws.Cells.Item(672, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(682, 9).Value = 'sv'
$ws.Cells.Item(682, 10).Value = 'Statement-opinion'
$ws.Cells.Item(698, 9).Value = 'sv'
$ws.Cells.Item(698, 10).Value = 'Statement-opinion'
$ws.Cells.Item(705, 9).Value = 'sv'
$ws.Cells.Item(705, 10).Value = 'Statement-opinion'
